$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 described "My Account" / "Account details" test case.
# Rename it to the ProductListPage / PLP details test case, since all
# test cases relate to the PLP page now.
$ws.Range("A3").Value = "ProductListPage"
$ws.Range("B3").Value = "PLP details"

# Update the active selection to B3 (matches the recorded view state).
$ws.Range("B3").Select()
